$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the statistics values for row 51 (copy counting)
$ws.Range("B51").Value = 2625.0014253999998
$ws.Range("C51").Value = 24847.628011699937
$ws.Range("D51").Value = 3384
$ws.Range("E51").Value = 8678.7405399999989

# Update the statistics values for row 95 (copy counting)
$ws.Range("B95").Value = 3058.1267253999999
$ws.Range("C95").Value = 28986.304910366598
$ws.Range("D95").Value = 3384
$ws.Range("E95").Value = 9268.3685399999995

# Update the statistics values for row 104 (process areas)
$ws.Range("B104").Value = 324.12076519999999
$ws.Range("C104").Value = 988.09399305185264
$ws.Range("D104").Value = 577
$ws.Range("E104").Value = 768.72852

# Update the view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 29
$ws.Range("A47").Select()
